# Aula 11 - "Divisão e Conquista" title update
# Slide 2 (section divider): "Aula 10 / Contextualização" -> "Aula 11 / Divisão e Conquista"
# Slide 3 (content title): "Contextualização" -> "Divisão e Conquista" (typed as 2 edits,
#   which is why it lands in the saved file as 3 runs: "Divisão" + " " + "e Conquista")

$p = $ppt.ActivePresentation

# --- Slide 2 -------------------------------------------------------------
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(3)
$tr2 = $sh2.TextFrame.TextRange

$aula = $tr2.Find("Aula 10", 0)
$aula.Text = "Aula 11"

$ctx2 = $tr2.Find("Contextualização", 0)
$ctx2.Text = "Divisão e Conquista"

# --- Slide 3 -------------------------------------------------------------
$s3  = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

# Shrink the existing run down to "Divisão" (keeps its original run formatting:
# language, bold, color, ...).
$ctx3 = $tr3.Find("Contextualização", 0)
$ctx3.Text = "Divisão"

# Append the rest of the title after that run, then re-assign the new substrings'
# text so the engine splits them into their own runs (" " and "e Conquista"),
# inheriting the same bold/color formatting as the rest of the title.
$null = $ctx3.InsertAfter(" e Conquista")

$sp3 = $tr3.Characters($ctx3.Start + $ctx3.Length, 1)
$sp3.Text = " "

$rest3 = $tr3.Characters($sp3.Start + 1, 11)
$rest3.Text = "e Conquista"
